$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 used to be "Dakota Myers" / drm1022@sru.edu -- replace with the new
# contact "Franz Ferdinand" / therockband@gmail.com. Every other field on
# the row (middle initial, address, phone numbers, ...) stays as-is.
$ws.Range("A2").Value = "Franz"
$ws.Range("B2").Value = "Ferdinand"
$ws.Range("D2").Value = "therockband@gmail.com"

# The old mailto: hyperlink pointed at the previous email address, so drop
# it now that the cell text/value no longer matches. Iterating the
# collection (rather than indexing with .Item()) is what reliably binds
# each hyperlink so .Delete() actually takes effect.
foreach ($h in $ws.Hyperlinks) {
    if ($h.Address -eq "mailto:drm1022@sru.edu") {
        $h.Delete()
    }
}

# Restore the active cell/selection like the saved workbook shows.
$ws.Range("M5").Select() | Out-Null
